# Update countries & provincias Spain
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Pais")

# --- Update the "last updated" timestamp text (row 1) ---
$ws.Range("A1").Value = "Datos actualizados a 29 de Marzo de 2020 a las 10:25"

# --- Reorder Rusia: it now appears right after Japon (row 33) and before
#     Pakistan. Pakistan, Rumania and Tailandia shift down one row, keeping
#     their own data values. Rusia gets refreshed totals. ---
$ws.Range("A34").Value = "Rusia"
$ws.Range("B34").Value = 1534
$ws.Range("C34").Value = 270
$ws.Range("D34").Value = 64
$ws.Range("E34").Value = 1462
$ws.Range("F34").Value = 8
$ws.Range("G34").Value = 4
$ws.Range("H34").Value = 8

$ws.Range("A35").Value = "Pakistan"
$ws.Range("B35").Value = 1526
$ws.Range("C35").Value = 31
$ws.Range("D35").Value = 29
$ws.Range("E35").Value = 1484
$ws.Range("F35").Value = 11
$ws.Range("G35").Value = 1
$ws.Range("H35").Value = 13

$ws.Range("A36").Value = "Rumania"
$ws.Range("B36").Value = 1452
$ws.Range("C36").Value = 0
$ws.Range("D36").Value = 139
$ws.Range("E36").Value = 1276
$ws.Range("F36").Value = 34
$ws.Range("G36").Value = 0
$ws.Range("H36").Value = 37

$ws.Range("A37").Value = "Tailandia"
$ws.Range("B37").Value = 1388
$ws.Range("C37").Value = 143
$ws.Range("D37").Value = 97
$ws.Range("E37").Value = 1284
$ws.Range("F37").Value = 11
$ws.Range("G37").Value = 1
$ws.Range("H37").Value = 7

# --- Other updated case counts ---

# Austria (row 16)
$ws.Range("B16").Value = 8346
$ws.Range("C16").Value = 75
$ws.Range("E16").Value = 8053

# Noruega (row 20)
$ws.Range("B20").Value = 4037
$ws.Range("C20").Value = 22
$ws.Range("E20").Value = 4007

# Chequia (row 25)
$ws.Range("E25").Value = 2639
$ws.Range("G25").Value = 2
$ws.Range("H25").Value = 13

# Albania (row 88)
$ws.Range("D88").Value = 33
$ws.Range("E88").Value = 154

# Sri Lanka (row 102)
$ws.Range("D102").Value = 10
$ws.Range("E102").Value = 104
